$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "94.908.87"
$ws.Range("E2").Value = "  -1.82%  "

$ws.Range("D3").Value = "3.482.22"
$ws.Range("E3").Value = "  +4.41%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'239.37"
$ws.Range("E5").Value = "  -4.16%  "

$ws.Range("D6").Value = "'643.96"
$ws.Range("E6").Value = "  -1.77%  "

$ws.Range("E7").Value = "  +3.38%  "

$ws.Range("D8").Value = "'0.402"
$ws.Range("E8").Value = "  -5.13%  "

$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("D10").Value = "'0.994"
$ws.Range("E10").Value = "  -0.86%  "

$ws.Range("D11").Value = "3.476.02"
$ws.Range("E11").Value = "  +4.31%  "

$ws.Range("D12").Value = "'42.69"
$ws.Range("E12").Value = "  +5.10%  "

$ws.Range("E13").Value = "  -3.73%  "

$ws.Range("E14").Value = "  +2.13%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "4.132.09"
$ws.Range("E15").Value = "  +4.34%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "94.753.29"
$ws.Range("E16").Value = "  -1.73%  "

$ws.Range("E17").Value = "  +2.17%  "

$ws.Range("D18").Value = "'8.46"
$ws.Range("E18").Value = "  -2.55%  "

$ws.Range("D19").Value = "3.493.96"
$ws.Range("E19").Value = "  +4.27%  "

$ws.Range("D20").Value = "'17.89"
$ws.Range("E20").Value = "  +2.50%  "

$ws.Range("E21").Value = "  +7.63%  "

$ws.Range("E22").Value = "  -7.92%  "

$ws.Range("D23").Value = "'505.09"
$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("D24").Value = "'3.19"
$ws.Range("E24").Value = "  -4.66%  "

$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "'0.0000192"
$ws.Range("E25").Value = "  -3.03%  "

$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").Value = "'6.57"
$ws.Range("E26").Value = "  -0.43%  "

$ws.Range("D27").Value = "'92.09"
$ws.Range("E27").Value = "  -4.40%  "

$ws.Range("D28").Value = "'12.10"
$ws.Range("E28").Value = "  -0.18%  "

$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'11.88"
$ws.Range("E29").Value = "  +5.66%  "

$ws.Range("B30").Value = "Dai"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  -0.19%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'2.77"
$ws.Range("E31").Value = "  +9.72%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.137"
$ws.Range("E32").Value = "  -5.11%  "

$ws.Range("B33").Value = "Cronos"
$ws.Range("C33").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D33").Value = "'0.184"
$ws.Range("E33").Value = "  -2.28%  "

$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").Value = "'0.997"
$ws.Range("E34").Value = "  -0.16%  "

$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "'30.58"
$ws.Range("E35").Value = "  +8.07%  "

$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D36").Value = "'0.575"
$ws.Range("E36").Value = "  +3.87%  "

$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").Value = "'547.62"
$ws.Range("E37").Value = "  +8.34%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").Value = "'7.75"
$ws.Range("E38").Value = "  -0.80%  "

$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").Value = "'1.46"
$ws.Range("E39").Value = "  -2.43%  "

$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "'0.945"
$ws.Range("E40").Value = "  +13.59%  "

$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  -0.02%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.151"
$ws.Range("E42").Value = "  +0.37%  "

$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").Value = "'24.09"
$ws.Range("E43").Value = "  -1.09%  "

$ws.Range("B44").Value = "ImmutableX"
$ws.Range("C44").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D44").Value = "'1.70"
$ws.Range("E44").Value = "  +1.82%  "

$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").Value = "'5.64"
$ws.Range("E45").Value = "  +1.45%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0415"
$ws.Range("E46").Value = "  -3.95%  "

$ws.Range("B47").Value = "MantraDAO"
$ws.Range("C47").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D47").Value = "'3.53"
$ws.Range("E47").Value = "  -3.19%  "

$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "'2.18"
$ws.Range("E48").Value = "  +10.09%  "

$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").Value = "'3.23"
$ws.Range("E49").Value = "  +3.74%  "

$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").Value = "'53.48"
$ws.Range("E50").Value = "  +0.06%  "

$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "'8.09"
$ws.Range("E51").Value = "  -4.35%  "
